$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6").Value = "Крупы"
[void]$ws.Range("A7").Select()
